$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = [double]"1.988074333333333"
$ws.Cells.Item(2, 8).Value = [double]"5.964223"
$ws.Cells.Item(2, 9).Value = [double]"0.01657769708907969"
$ws.Cells.Item(2, 10).Value = [double]"0.01657769708907968"
$ws.Cells.Item(2, 13).Value = [double]"1.918906333333333"
$ws.Cells.Item(2, 14).Value = [double]"5.756718999999999"
$ws.Cells.Item(2, 15).Value = [double]"0.006524019162508824"
$ws.Cells.Item(2, 16).Value = [double]"0.006524019162508824"
$ws.Cells.Item(2, 17).Value = [double]"3.814928429370777"
$ws.Cells.Item(2, 18).Value = [double]"34.33435586433699"
$ws.Cells.Item(2, 19).Value = [double]"0.0001081532134794226"
$ws.Cells.Item(2, 20).Value = [double]"0.0001081532134794226"
$ws.Cells.Item(3, 7).Value = [double]"1.988074333333333"
$ws.Cells.Item(3, 8).Value = [double]"5.964223"
$ws.Cells.Item(3, 9).Value = [double]"0.01657769708907969"
$ws.Cells.Item(3, 10).Value = [double]"0.01657769708907968"
$ws.Cells.Item(3, 15).Value = [double]"0.6163557430885885"
$ws.Cells.Item(3, 16).Value = [double]"0.6163557430885885"
$ws.Cells.Item(3, 17).Value = [double]"360.4147977410892"
$ws.Cells.Item(3, 18).Value = [double]"3243.733179669802"
$ws.Cells.Item(3, 19).Value = [double]"0.01021775880803724"
$ws.Cells.Item(3, 20).Value = [double]"0.01021775880803724"
$ws.Cells.Item(4, 7).Value = [double]"1.988074333333333"
$ws.Cells.Item(4, 8).Value = [double]"5.964223"
$ws.Cells.Item(4, 9).Value = [double]"0.01657769708907969"
$ws.Cells.Item(4, 10).Value = [double]"0.01657769708907968"
$ws.Cells.Item(4, 13).Value = [double]"29.04767233333333"
$ws.Cells.Item(4, 14).Value = [double]"87.143017"
$ws.Cells.Item(4, 15).Value = [double]"0.09875811426384234"
$ws.Cells.Item(4, 16).Value = [double]"0.09875811426384236"
$ws.Cells.Item(4, 17).Value = [double]"57.74893180897678"
$ws.Cells.Item(4, 18).Value = [double]"519.740386280791"
$ws.Cells.Item(4, 19).Value = [double]"0.001637182103354698"
$ws.Cells.Item(4, 20).Value = [double]"0.001637182103354698"
$ws.Cells.Item(5, 7).Value = [double]"1.988074333333333"
$ws.Cells.Item(5, 8).Value = [double]"5.964223"
$ws.Cells.Item(5, 9).Value = [double]"0.01657769708907969"
$ws.Cells.Item(5, 10).Value = [double]"0.01657769708907968"
$ws.Cells.Item(5, 13).Value = [double]"81.87450533333333"
$ws.Cells.Item(5, 14).Value = [double]"245.623516"
$ws.Cells.Item(5, 15).Value = [double]"0.2783621234850603"
$ws.Cells.Item(5, 16).Value = [double]"0.2783621234850603"
$ws.Cells.Item(5, 17).Value = [double]"162.7726026075631"
$ws.Cells.Item(5, 18).Value = [double]"1464.953423468068"
$ws.Cells.Item(5, 19).Value = [double]"0.004614602964208324"
$ws.Cells.Item(5, 20).Value = [double]"0.004614602964208323"
$ws.Cells.Item(6, 9).Value = [double]"0.7746030815641455"
$ws.Cells.Item(6, 10).Value = [double]"0.7746030815641454"
$ws.Cells.Item(6, 13).Value = [double]"1.918906333333333"
$ws.Cells.Item(6, 14).Value = [double]"5.756718999999999"
$ws.Cells.Item(6, 15).Value = [double]"0.006524019162508824"
$ws.Cells.Item(6, 16).Value = [double]"0.006524019162508824"
$ws.Cells.Item(6, 17).Value = [double]"178.2548746944995"
$ws.Cells.Item(6, 18).Value = [double]"1604.293872250496"
$ws.Cells.Item(6, 19).Value = [double]"0.005053525347462871"
$ws.Cells.Item(6, 20).Value = [double]"0.00505352534746287"
$ws.Cells.Item(7, 9).Value = [double]"0.7746030815641455"
$ws.Cells.Item(7, 10).Value = [double]"0.7746030815641454"
$ws.Cells.Item(7, 15).Value = [double]"0.6163557430885885"
$ws.Cells.Item(7, 16).Value = [double]"0.6163557430885885"
$ws.Cells.Item(7, 19).Value = [double]"0.4774310579361795"
$ws.Cells.Item(7, 20).Value = [double]"0.4774310579361793"
$ws.Cells.Item(8, 9).Value = [double]"0.7746030815641455"
$ws.Cells.Item(8, 10).Value = [double]"0.7746030815641454"
$ws.Cells.Item(8, 13).Value = [double]"29.04767233333333"
$ws.Cells.Item(8, 14).Value = [double]"87.143017"
$ws.Cells.Item(8, 15).Value = [double]"0.09875811426384234"
$ws.Cells.Item(8, 16).Value = [double]"0.09875811426384236"
$ws.Cells.Item(8, 17).Value = [double]"2698.354318811748"
$ws.Cells.Item(8, 18).Value = [double]"24285.18886930573"
$ws.Cells.Item(8, 19).Value = [double]"0.07649833963823627"
$ws.Cells.Item(8, 20).Value = [double]"0.07649833963823627"
$ws.Cells.Item(9, 9).Value = [double]"0.7746030815641455"
$ws.Cells.Item(9, 10).Value = [double]"0.7746030815641454"
$ws.Cells.Item(9, 13).Value = [double]"81.87450533333333"
$ws.Cells.Item(9, 14).Value = [double]"245.623516"
$ws.Cells.Item(9, 15).Value = [double]"0.2783621234850603"
$ws.Cells.Item(9, 16).Value = [double]"0.2783621234850603"
$ws.Cells.Item(9, 17).Value = [double]"7605.649861770638"
$ws.Cells.Item(9, 18).Value = [double]"68450.84875593575"
$ws.Cells.Item(9, 19).Value = [double]"0.2156201586422669"
$ws.Cells.Item(9, 20).Value = [double]"0.2156201586422669"
$ws.Cells.Item(10, 7).Value = [double]"23.741365"
$ws.Cells.Item(10, 8).Value = [double]"71.22409500000001"
$ws.Cells.Item(10, 9).Value = [double]"0.1979690350870239"
$ws.Cells.Item(10, 10).Value = [double]"0.1979690350870239"
$ws.Cells.Item(10, 13).Value = [double]"1.918906333333333"
$ws.Cells.Item(10, 14).Value = [double]"5.756718999999999"
$ws.Cells.Item(10, 15).Value = [double]"0.006524019162508824"
$ws.Cells.Item(10, 16).Value = [double]"0.006524019162508824"
$ws.Cells.Item(10, 17).Value = [double]"45.55745566047833"
$ws.Cells.Item(10, 18).Value = [double]"410.017100944305"
$ws.Cells.Item(10, 19).Value = [double]"0.001291553778491126"
$ws.Cells.Item(10, 20).Value = [double]"0.001291553778491126"
$ws.Cells.Item(11, 7).Value = [double]"23.741365"
$ws.Cells.Item(11, 8).Value = [double]"71.22409500000001"
$ws.Cells.Item(11, 9).Value = [double]"0.1979690350870239"
$ws.Cells.Item(11, 10).Value = [double]"0.1979690350870239"
$ws.Cells.Item(11, 15).Value = [double]"0.6163557430885885"
$ws.Cells.Item(11, 16).Value = [double]"0.6163557430885885"
$ws.Cells.Item(11, 17).Value = [double]"4304.033868907504"
$ws.Cells.Item(11, 18).Value = [double]"38736.30482016753"
$ws.Cells.Item(11, 19).Value = [double]"0.1220193517295935"
$ws.Cells.Item(11, 20).Value = [double]"0.1220193517295935"
$ws.Cells.Item(12, 7).Value = [double]"23.741365"
$ws.Cells.Item(12, 8).Value = [double]"71.22409500000001"
$ws.Cells.Item(12, 9).Value = [double]"0.1979690350870239"
$ws.Cells.Item(12, 10).Value = [double]"0.1979690350870239"
$ws.Cells.Item(12, 13).Value = [double]"29.04767233333333"
$ws.Cells.Item(12, 14).Value = [double]"87.143017"
$ws.Cells.Item(12, 15).Value = [double]"0.09875811426384234"
$ws.Cells.Item(12, 16).Value = [double]"0.09875811426384236"
$ws.Cells.Item(12, 17).Value = [double]"689.6313912660684"
$ws.Cells.Item(12, 18).Value = [double]"6206.682521394616"
$ws.Cells.Item(12, 19).Value = [double]"0.01955104858782692"
$ws.Cells.Item(12, 20).Value = [double]"0.01955104858782693"
$ws.Cells.Item(13, 7).Value = [double]"23.741365"
$ws.Cells.Item(13, 8).Value = [double]"71.22409500000001"
$ws.Cells.Item(13, 9).Value = [double]"0.1979690350870239"
$ws.Cells.Item(13, 10).Value = [double]"0.1979690350870239"
$ws.Cells.Item(13, 13).Value = [double]"81.87450533333333"
$ws.Cells.Item(13, 14).Value = [double]"245.623516"
$ws.Cells.Item(13, 15).Value = [double]"0.2783621234850603"
$ws.Cells.Item(13, 16).Value = [double]"0.2783621234850603"
$ws.Cells.Item(13, 17).Value = [double]"1943.812515313113"
$ws.Cells.Item(13, 18).Value = [double]"17494.31263781802"
$ws.Cells.Item(13, 19).Value = [double]"0.05510708099111239"
$ws.Cells.Item(13, 20).Value = [double]"0.05510708099111239"
$ws.Cells.Item(14, 7).Value = [double]"1.301204666666667"
$ws.Cells.Item(14, 8).Value = [double]"3.903614"
$ws.Cells.Item(14, 9).Value = [double]"0.01085018625975097"
$ws.Cells.Item(14, 10).Value = [double]"0.01085018625975097"
$ws.Cells.Item(14, 13).Value = [double]"1.918906333333333"
$ws.Cells.Item(14, 14).Value = [double]"5.756718999999999"
$ws.Cells.Item(14, 15).Value = [double]"0.006524019162508824"
$ws.Cells.Item(14, 16).Value = [double]"0.006524019162508824"
$ws.Cells.Item(14, 17).Value = [double]"2.496889875829555"
$ws.Cells.Item(14, 18).Value = [double]"22.472008882466"
$ws.Cells.Item(14, 19).Value = [double]"7.078682307540529E-05"
$ws.Cells.Item(14, 20).Value = [double]"7.078682307540527E-05"
$ws.Cells.Item(15, 7).Value = [double]"1.301204666666667"
$ws.Cells.Item(15, 8).Value = [double]"3.903614"
$ws.Cells.Item(15, 9).Value = [double]"0.01085018625975097"
$ws.Cells.Item(15, 10).Value = [double]"0.01085018625975097"
$ws.Cells.Item(15, 15).Value = [double]"0.6163557430885885"
$ws.Cells.Item(15, 16).Value = [double]"0.6163557430885885"
$ws.Cells.Item(15, 17).Value = [double]"235.8933008154263"
$ws.Cells.Item(15, 18).Value = [double]"2123.039707338836"
$ws.Cells.Item(15, 19).Value = [double]"0.006687574614778403"
$ws.Cells.Item(15, 20).Value = [double]"0.006687574614778403"
$ws.Cells.Item(16, 7).Value = [double]"1.301204666666667"
$ws.Cells.Item(16, 8).Value = [double]"3.903614"
$ws.Cells.Item(16, 9).Value = [double]"0.01085018625975097"
$ws.Cells.Item(16, 10).Value = [double]"0.01085018625975097"
$ws.Cells.Item(16, 13).Value = [double]"29.04767233333333"
$ws.Cells.Item(16, 14).Value = [double]"87.143017"
$ws.Cells.Item(16, 15).Value = [double]"0.09875811426384234"
$ws.Cells.Item(16, 16).Value = [double]"0.09875811426384236"
$ws.Cells.Item(16, 17).Value = [double]"37.79696679593756"
$ws.Cells.Item(16, 18).Value = [double]"340.172701163438"
$ws.Cells.Item(16, 19).Value = [double]"0.001071543934424459"
$ws.Cells.Item(16, 20).Value = [double]"0.001071543934424459"
$ws.Cells.Item(17, 7).Value = [double]"1.301204666666667"
$ws.Cells.Item(17, 8).Value = [double]"3.903614"
$ws.Cells.Item(17, 9).Value = [double]"0.01085018625975097"
$ws.Cells.Item(17, 10).Value = [double]"0.01085018625975097"
$ws.Cells.Item(17, 13).Value = [double]"81.87450533333333"
$ws.Cells.Item(17, 14).Value = [double]"245.623516"
$ws.Cells.Item(17, 15).Value = [double]"0.2783621234850603"
$ws.Cells.Item(17, 16).Value = [double]"0.2783621234850603"
$ws.Cells.Item(17, 17).Value = [double]"106.5354884207582"
$ws.Cells.Item(17, 18).Value = [double]"958.819395786824"
$ws.Cells.Item(17, 19).Value = [double]"0.003020280887472704"
$ws.Cells.Item(17, 20).Value = [double]"0.003020280887472704"
